# Update the "Estado de Cuenta" (EC) database rows on Hoja1.
# The late-payment detail table occupies rows 16-20 (cols C:F):
#   C = worker document number, D = worker name, E = late period, F = overdue amount.
# This edit re-sorts/refreshes that data: a new record for EDELMIRA CASTILLO ALVAREZ
# (period 2303) is folded in alongside KETTY LUZ ACOSTA MARTINEZ's two periods, and
# NORELIS MENDOZA ROCHA's two periods swap order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 keeps KETTY LUZ ACOSTA MARTINEZ, but now reflects period 2303.
$ws.Range("C16").Value = "1049452250"
$ws.Range("D16").Value = "KETTY LUZ ACOSTA MARTINEZ"
$ws.Range("E16").Value = "2303"
$ws.Range("F16").Value = 46400

# Row 17 now holds EDELMIRA CASTILLO ALVAREZ, period 2303.
$ws.Range("C17").Value = "45372092"
$ws.Range("D17").Value = "EDELMIRA CASTILLO ALVAREZ"
$ws.Range("E17").Value = "2303"
$ws.Range("F17").Value = 13920

# Row 18 now holds KETTY LUZ ACOSTA MARTINEZ, period 2304.
$ws.Range("C18").Value = "1049452250"
$ws.Range("D18").Value = "KETTY LUZ ACOSTA MARTINEZ"
$ws.Range("E18").Value = "2304"
$ws.Range("F18").Value = 46400

# Rows 19 & 20 (NORELIS MENDOZA ROCHA) swap their period/amount pairing.
$ws.Range("E19").Value = "2307"
$ws.Range("F19").Value = 46400

$ws.Range("E20").Value = "2308"
$ws.Range("F20").Value = 6187
